$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing emails for Marry Smith / John Smith (row 3 / row 2)
$ws.Range("C3").Value2 = "g369keh0@freeml.net"
$ws.Range("C2").Value2 = "rndgmvkcf@emlpro.com"

# Add a new row of data (row 4): Kira Alex, Russia, karzanovalexey@gmail.com
$ws.Range("B4").Value2 = "Kira"
$ws.Range("A4").Value2 = "Alex"
$ws.Range("D4").Value2 = "Russia"
$ws.Range("C4").Value2 = "karzanovalexey@gmail.com"

# Turn the new email address into a mailto hyperlink
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:karzanovalexey@gmail.com")

# Move the active selection to E10 (matches author's final cursor position)
$ws.Range("E10").Select()
